# on-call-rotation.xlsx tutorial refresh:
#  - roll every "On Call Start/End Date" forward by exactly two years
#    (the tutorial sample data is periodically bumped so it still looks
#    "current" for readers), and
#  - leave the selection on K11 (last place the author clicked before
#    saving).
#
# NOTE: the source diff also touches a handful of attributes that are
# purely save-environment artifacts of the Excel client/tenant that
# produced the file -- fileVersion/rupBuild, xr:revisionPtr coauthoring
# GUIDs, the x15ac:absPath (a OneDrive/SharePoint URL baked in by that
# user's machine), workbookView window-geometry (xWindow/yWindow/
# windowWidth/windowHeight), the theme's internal display name ("Office"
# -> "Office 2013 - 2022"), sheetFormatPr's defaultRowHeight/dyDescent,
# and sub-pixel "best fit" column-width rounding, plus the
# docMetadata/LabelInfo.xml MIP sensitivity-label part that a SharePoint
# tenant stamps onto files it hosts. None of these are exposed as
# settable state on the Excel object model (real or emulated) -- they
# are written by the host application/service itself at save time, not
# by user/script action -- so they are intentionally left alone here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("H1")
$ws.Activate()

$dates = @(
    @{ Cell = "C3";  Value = 45295 }
    @{ Cell = "D3";  Value = 45302 }
    @{ Cell = "C4";  Value = 45302 }
    @{ Cell = "D4";  Value = 45309 }
    @{ Cell = "C5";  Value = 45309 }
    @{ Cell = "D5";  Value = 45316 }
    @{ Cell = "C6";  Value = 45316 }
    @{ Cell = "D6";  Value = 45323 }
    @{ Cell = "C7";  Value = 45323 }
    @{ Cell = "D7";  Value = 45330 }
    @{ Cell = "C8";  Value = 45330 }
    @{ Cell = "D8";  Value = 45337 }
    @{ Cell = "C9";  Value = 45337 }
    @{ Cell = "D9";  Value = 45344 }
    @{ Cell = "C10"; Value = 45344 }
    @{ Cell = "D10"; Value = 45352 }
    @{ Cell = "C11"; Value = 45352 }
    @{ Cell = "D11"; Value = 45359 }
    @{ Cell = "C12"; Value = 45359 }
    @{ Cell = "D12"; Value = 45366 }
    @{ Cell = "C13"; Value = 45366 }
    @{ Cell = "D13"; Value = 45373 }
    @{ Cell = "C14"; Value = 45373 }
    @{ Cell = "D14"; Value = 45380 }
    @{ Cell = "C15"; Value = 45380 }
    @{ Cell = "D15"; Value = 45387 }
    @{ Cell = "C16"; Value = 45387 }
    @{ Cell = "D16"; Value = 45394 }
    @{ Cell = "C17"; Value = 45394 }
    @{ Cell = "D17"; Value = 45401 }
    @{ Cell = "C18"; Value = 45401 }
    @{ Cell = "D18"; Value = 45408 }
    @{ Cell = "C19"; Value = 45408 }
    @{ Cell = "D19"; Value = 45415 }
    @{ Cell = "C20"; Value = 45415 }
    @{ Cell = "D20"; Value = 45422 }
    @{ Cell = "C21"; Value = 45422 }
    @{ Cell = "D21"; Value = 45429 }
    @{ Cell = "C22"; Value = 45429 }
    @{ Cell = "D22"; Value = 45437 }
    @{ Cell = "C23"; Value = 45437 }
    @{ Cell = "D23"; Value = 45302 }
    @{ Cell = "C24"; Value = 45443 }
    @{ Cell = "D24"; Value = 45309 }
    @{ Cell = "C25"; Value = 45450 }
    @{ Cell = "D25"; Value = 45457 }
    @{ Cell = "C26"; Value = 45457 }
    @{ Cell = "D26"; Value = 45464 }
    @{ Cell = "C27"; Value = 45464 }
    @{ Cell = "D27"; Value = 45471 }
    @{ Cell = "C28"; Value = 45471 }
    @{ Cell = "D28"; Value = 45478 }
)

foreach ($d in $dates) {
    $ws.Range($d.Cell).Value = $d.Value
}

# Match the author's final selection before the file was saved.
$ws.Range("K11").Select()
